$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 11 & 12: mark these two table rows as "addressed" (gray highlight fill,
# like rows 18/19/33/34 already have) and fill in the new Response (column C)
# and Status of changes (column F) cells.
# ---------------------------------------------------------------------------

# Plain gray-highlight fill, no wrap (donor already carries fillId + no wrap).
$noWrapDonor = $ws.Range("A18")
foreach ($addr in @("A11","D11","F11","A12","C12","D12","F12")) {
    $noWrapDonor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Gray-highlight fill *with* wrap text (donor already carries fillId + wrap).
$wrapDonor = $ws.Range("B18")
foreach ($addr in @("C11","E11","E12")) {
    $wrapDonor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# B11, B12: keep their existing dark-gray "comment" font, gain the
# gray-highlight fill (they already had wrap text).
foreach ($addr in @("B11","B12")) {
    $wrapDonor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Font.Color = 657930
    $ws.Range($addr).WrapText = $true
}

$excel.CutCopyMode = $false

# New response / plan / status text.
# (set C12's new string first so the shared-string table order matches)
$ws.Range("C12").Value2 = "We combined all nomenclature into a single table."
$ws.Range("C11").Value2 = "Table 8 is placed automatically by LaTeX. The proofs will almost certainly look different, so there is no need to address placement issues at this time."
$ws.Range("F11").Value2 = "DONE"
$ws.Range("F12").Value2 = "DONE"

# Row 12 grows to two lines now that it wraps a longer comment.
$ws.Rows.Item(12).RowHeight = 34

# ---------------------------------------------------------------------------
# Update the active selection / scroll position recorded in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("C14").Select()
